$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the last existing row (366) down into the new rows (367-374)
$ws.Range("A366:D366").Copy($ws.Range("A367:D374"))

$ws.Range("A367").Value = 44441
$ws.Range("B367").Value = 3
$ws.Range("C367").Value = 10
$ws.Range("D367").Value = 29.0985276145027

$ws.Range("A368").Value = 44442
$ws.Range("B368").Value = 2
$ws.Range("C368").Value = 8
$ws.Range("D368").Value = 23.27882209160216

$ws.Range("A369").Value = 44443
$ws.Range("B369").Value = 1
$ws.Range("C369").Value = 9
$ws.Range("D369").Value = 26.18867485305244

$ws.Range("A370").Value = 44444
$ws.Range("B370").Value = 1
$ws.Range("C370").Value = 7
$ws.Range("D370").Value = 20.36896933015189

$ws.Range("A371").Value = 44445
$ws.Range("B371").Value = 1
$ws.Range("C371").Value = 8
$ws.Range("D371").Value = 23.27882209160216

$ws.Range("A372").Value = 44446
$ws.Range("B372").Value = 0
$ws.Range("C372").Value = 8
$ws.Range("D372").Value = 23.27882209160216

$ws.Range("A373").Value = 44447
$ws.Range("B373").Value = 0
$ws.Range("C373").Value = 8
$ws.Range("D373").Value = 23.27882209160216

$ws.Range("A374").Value = 44448
$ws.Range("B374").Value = 1
$ws.Range("C374").Value = 6
$ws.Range("D374").Value = 17.45911656870162
